$d = $word.ActiveDocument

$ids = @("p042r_1", "p042r_2", "p042r_3")
foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
